# Apply small adjustments to DPE for EPICP to reflect dataset file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F6: "PAL" -> "pal"
$ws.Range("F6").Value = "pal"

# Populate column F (input_variables) for rows 10-128 based on existing
# column B (dataschema_variable) and column G (rule_category) values:
#   - if rule_category ("G") is "impossible", input_variables is "impossible"
#   - otherwise, input_variables takes the same value as dataschema_variable ("B")
for ($row = 10; $row -le 128; $row++) {
    $ruleCategory = $ws.Cells.Item($row, 7).Value2
    if ($ruleCategory -eq "impossible") {
        $ws.Cells.Item($row, 6).Value = "impossible"
    } else {
        $varName = $ws.Cells.Item($row, 2).Value2
        $ws.Cells.Item($row, 6).Value = $varName
    }
}
